$p = $ppt.ActivePresentation
try {
  $p.ApplyTheme("C:\nonexistent\theme1.thmx")
  Write-Host "no error"
} catch {
  Write-Host "ERR: $_"
}
